# Actualización automática 2025-12-07 17:44:30
# Updates the PRESUPUESTO (column G) figures on the "VENTA MENSUAL" sheet
# for the listed clients, then refreshes the column total in G62.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Row number -> new PRESUPUESTO (column G) value
$updates = @{
    5  = 0
    6  = 500
    8  = 2000
    15 = 5000
    22 = 2000
    25 = 5000
    27 = 2000
    28 = 0
    29 = 1500
    30 = 500
    31 = 0
    32 = 5000
    33 = 8000
    34 = 5000
    35 = 8000
    40 = 1000
    42 = 0
    43 = 5000
    44 = 1500
    48 = 1500
    49 = 2000
    50 = 3000
    51 = 5000
    52 = 1500
    54 = 5000
    58 = 0
    59 = 2000
    60 = 500
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}

# Recompute the totals row (row 62) for column G from rows 2-61, matching
# the workbook's existing pattern of storing the sum as a literal value.
$total = 0
for ($r = 2; $r -le 61; $r++) {
    $total += $ws.Cells.Item($r, 7).Value2
}
$ws.Cells.Item(62, 7).Value = $total
